# Update "loading_percent" values for the 380 kV case (Case_4_158)
# Rows 2-25 correspond to time steps 0-23; columns B,C,D,E,G,H,J,K,M hold
# the recomputed line-loading percentages (F, I, L, N, O remain 0 and are untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.45616125726453
$ws.Range("C2").Value = 5.540912480444399
$ws.Range("D2").Value = 14.79634839485919
$ws.Range("E2").Value = 15.85121562336859
$ws.Range("G2").Value = 68.99953833991569
$ws.Range("H2").Value = 23.93706641619157
$ws.Range("J2").Value = 9.266839676897789
$ws.Range("K2").Value = 13.09487274637254
$ws.Range("M2").Value = 18.63421141517745

$ws.Range("B3").Value = 13.37336625490839
$ws.Range("C3").Value = 5.518470734537807
$ws.Range("D3").Value = 14.77171769822337
$ws.Range("E3").Value = 15.84451855489243
$ws.Range("G3").Value = 68.29266798032268
$ws.Range("H3").Value = 23.84782489767796
$ws.Range("J3").Value = 9.281490311776777
$ws.Range("K3").Value = 13.07119969992558
$ws.Range("M3").Value = 18.64144191972745

$ws.Range("B4").Value = 13.32709030605761
$ws.Range("C4").Value = 5.507552550882399
$ws.Range("D4").Value = 14.75949827461947
$ws.Range("E4").Value = 15.84310253950619
$ws.Range("G4").Value = 67.86301473424442
$ws.Range("H4").Value = 23.79519013036274
$ws.Range("J4").Value = 9.291301483799904
$ws.Range("K4").Value = 13.06092410767734
$ws.Range("M4").Value = 18.65005450795602

$ws.Range("B5").Value = 13.30939920030442
$ws.Range("C5").Value = 5.503827407929085
$ws.Range("D5").Value = 14.75525185073599
$ws.Range("E5").Value = 15.84320435253398
$ws.Range("G5").Value = 67.68915620627682
$ws.Range("H5").Value = 23.77429204857863
$ws.Range("J5").Value = 9.295504999718142
$ws.Range("K5").Value = 13.05781191608631
$ws.Range("M5").Value = 18.65461381670037

$ws.Range("B6").Value = 13.30653261225367
$ws.Range("C6").Value = 5.503252692776278
$ws.Range("D6").Value = 14.75459108524571
$ws.Range("E6").Value = 15.84326228004443
$ws.Range("G6").Value = 67.66036508305784
$ws.Range("H6").Value = 23.77085543685768
$ws.Range("J6").Value = 9.296215403449372
$ws.Range("K6").Value = 13.05736017334661
$ws.Range("M6").Value = 18.65543427865108

$ws.Range("B7").Value = 13.3268469689345
$ws.Range("C7").Value = 5.507499375370243
$ws.Range("D7").Value = 14.75943803422017
$ws.Range("E7").Value = 15.84310116297455
$ws.Range("G7").Value = 67.86066487292491
$ws.Range("H7").Value = 23.79490605060858
$ws.Range("J7").Value = 9.291357341916603
$ws.Range("K7").Value = 13.06087777771322
$ws.Range("M7").Value = 18.65011174659096

$ws.Range("B8").Value = 13.42667963274037
$ws.Range("C8").Value = 5.532583964648063
$ws.Range("D8").Value = 14.78725490742355
$ws.Range("E8").Value = 15.84834755476221
$ws.Range("G8").Value = 68.75498189594154
$ws.Range("H8").Value = 23.90585255305852
$ws.Range("J8").Value = 9.271722106887632
$ws.Range("K8").Value = 13.08582894938575
$ws.Range("M8").Value = 18.63583867813322

$ws.Range("B9").Value = 13.65765671441999
$ws.Range("C9").Value = 5.604202249376259
$ws.Range("D9").Value = 14.86469626931264
$ws.Range("E9").Value = 15.87997617760616
$ws.Range("G9").Value = 70.53755491831689
$ws.Range("H9").Value = 24.14018908193216
$ws.Range("J9").Value = 9.239676418926205
$ws.Range("K9").Value = 13.16832557769748
$ws.Range("M9").Value = 18.64093055645164

$ws.Range("B10").Value = 13.84736760491092
$ws.Range("C10").Value = 5.670027836403204
$ws.Range("D10").Value = 14.93531657239925
$ws.Range("E10").Value = 15.91613876371208
$ws.Range("G10").Value = 71.85685285226101
$ws.Range("H10").Value = 24.32206041227521
$ws.Range("J10").Value = 9.220053024056588
$ws.Range("K10").Value = 13.24901133874148
$ws.Range("M10").Value = 18.66477195489796

$ws.Range("B11").Value = 13.93767307266384
$ws.Range("C11").Value = 5.702715515693137
$ws.Range("D11").Value = 14.97036289672628
$ws.Range("E11").Value = 15.93537193834701
$ws.Range("G11").Value = 72.457426774773
$ws.Range("H11").Value = 24.40678424641884
$ws.Range("J11").Value = 9.211973547298792
$ws.Range("K11").Value = 13.28996975225819
$ws.Range("M11").Value = 18.6799579312883

$ws.Range("B12").Value = 13.97241439048144
$ws.Range("C12").Value = 5.715475721905318
$ws.Range("D12").Value = 14.98404797004936
$ws.Range("E12").Value = 15.94305246762041
$ws.Range("G12").Value = 72.68476559157756
$ws.Range("H12").Value = 24.4391418483442
$ws.Range("J12").Value = 9.209035610438203
$ws.Range("K12").Value = 13.3060806340529
$ws.Range("M12").Value = 18.68632958650821

$ws.Range("B13").Value = 13.96490851517733
$ws.Range("C13").Value = 5.712710797514855
$ws.Range("D13").Value = 14.98108234574457
$ws.Range("E13").Value = 15.94138070435162
$ws.Range("G13").Value = 72.63581021766714
$ws.Range("H13").Value = 24.43216103950208
$ws.Range("J13").Value = 9.209662944238861
$ws.Range("K13").Value = 13.30258433064321
$ws.Range("M13").Value = 18.68492976980563

$ws.Range("B14").Value = 13.94052053052513
$ws.Range("C14").Value = 5.703757722406832
$ws.Range("D14").Value = 14.97148052042774
$ws.Range("E14").Value = 15.93599587396511
$ws.Range("G14").Value = 72.47613261931824
$ws.Range("H14").Value = 24.409440891893
$ws.Range("J14").Value = 9.211729406177453
$ws.Range("K14").Value = 13.29128322910531
$ws.Range("M14").Value = 18.68046970391379

$ws.Range("B15").Value = 13.92565216003406
$ws.Range("C15").Value = 5.6983230799659
$ws.Range("D15").Value = 14.96565281951303
$ws.Range("E15").Value = 15.93274916700901
$ws.Range("G15").Value = 72.3783100438762
$ws.Range("H15").Value = 24.39555954506691
$ws.Range("J15").Value = 9.213011000210555
$ws.Range("K15").Value = 13.28443888078239
$ws.Range("M15").Value = 18.67781856373392

$ws.Range("B16").Value = 13.84154387004422
$ws.Range("C16").Value = 5.667945786046587
$ws.Range("D16").Value = 14.93308446727417
$ws.Range("E16").Value = 15.91493759290405
$ws.Range("G16").Value = 71.81759987684221
$ws.Range("H16").Value = 24.31656256755605
$ws.Range("J16").Value = 9.220598064964765
$ws.Range("K16").Value = 13.24641930435085
$ws.Range("M16").Value = 18.66386660378126

$ws.Range("B17").Value = 13.79095002524152
$ws.Range("C17").Value = 5.650004438025871
$ws.Range("D17").Value = 14.91384866974679
$ws.Range("E17").Value = 15.90472158597739
$ws.Range("G17").Value = 71.47362556577981
$ws.Range("H17").Value = 24.26860198320768
$ws.Range("J17").Value = 9.225469316850328
$ws.Range("K17").Value = 13.2241776695069
$ws.Range("M17").Value = 18.65641718164512

$ws.Range("B18").Value = 13.76222871094265
$ws.Range("C18").Value = 5.639944071875638
$ws.Range("D18").Value = 14.90306007300498
$ws.Range("E18").Value = 15.89910771241932
$ws.Range("G18").Value = 71.27582871288489
$ws.Range("H18").Value = 24.24120432288143
$ws.Range("J18").Value = 9.228350899749206
$ws.Range("K18").Value = 13.21178600365381
$ws.Range("M18").Value = 18.65254122946258

$ws.Range("B19").Value = 13.75257018406901
$ws.Range("C19").Value = 5.636582652657184
$ws.Range("D19").Value = 14.89945470757753
$ws.Range("E19").Value = 15.89725204745337
$ws.Range("G19").Value = 71.20887081189701
$ws.Range("H19").Value = 24.23196058027765
$ws.Range("J19").Value = 9.229340262136249
$ws.Range("K19").Value = 13.20765960856846
$ws.Range("M19").Value = 18.65129918863206

$ws.Range("B20").Value = 13.79629683229421
$ws.Range("C20").Value = 5.651887603841681
$ws.Range("D20").Value = 14.91586790148375
$ws.Range("E20").Value = 15.90578198985315
$ws.Range("G20").Value = 71.51023830248695
$ws.Range("H20").Value = 24.27368809469908
$ws.Range("J20").Value = 9.224942509985608
$ws.Range("K20").Value = 13.22650388937998
$ws.Range("M20").Value = 18.65716790117122

$ws.Range("B21").Value = 13.94766934758886
$ws.Range("C21").Value = 5.706377191475283
$ws.Range("D21").Value = 14.97428963124918
$ws.Range("E21").Value = 15.93756676933811
$ws.Range("G21").Value = 72.52303727544414
$ws.Range("H21").Value = 24.4161069940361
$ws.Range("J21").Value = 9.211119138458844
$ws.Range("K21").Value = 13.29458642264977
$ws.Range("M21").Value = 18.68176290537646

$ws.Range("B22").Value = 14.04975911887013
$ws.Range("C22").Value = 5.744209976477689
$ws.Range("D22").Value = 15.01488024941514
$ws.Range("E22").Value = 15.96065444897946
$ws.Range("G22").Value = 73.18439512034334
$ws.Range("H22").Value = 24.51078043192071
$ws.Range("J22").Value = 9.202793340605494
$ws.Range("K22").Value = 13.34257853984855
$ws.Range("M22").Value = 18.70145534030315

$ws.Range("B23").Value = 13.99499344787156
$ws.Range("C23").Value = 5.723819061602923
$ws.Range("D23").Value = 14.9929980699566
$ws.Range("E23").Value = 15.94812135410262
$ws.Range("G23").Value = 72.83151543501518
$ws.Range("H23").Value = 24.46010948641397
$ws.Range("J23").Value = 9.207172225824484
$ws.Range("K23").Value = 13.31664820975176
$ws.Range("M23").Value = 18.69061521203995

$ws.Range("B24").Value = 13.79387839915691
$ws.Range("C24").Value = 5.651035431364234
$ws.Range("D24").Value = 14.91495416435781
$ws.Range("E24").Value = 15.9053017729297
$ws.Range("G24").Value = 71.49368580214403
$ws.Range("H24").Value = 24.27138811603544
$ws.Range("J24").Value = 9.225180426935049
$ws.Range("K24").Value = 13.22545097309911
$ws.Range("M24").Value = 18.65682723345277

$ws.Range("B25").Value = 13.5915449407375
$ws.Range("C25").Value = 5.582469367332507
$ws.Range("D25").Value = 14.84131722236009
$ws.Range("E25").Value = 15.86914257271667
$ws.Range("G25").Value = 70.05305001034043
$ws.Range("H25").Value = 24.07505086768534
$ws.Range("J25").Value = 9.247655902156668
$ws.Range("K25").Value = 13.14244986666716
$ws.Range("M25").Value = 18.63601634044282
